$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws.Range("B2").Value = 0.1027396222441712
$ws.Range("B3").Value = 68.48362687895347
$ws.Range("B4").Value = 68.68568934491481
$ws.Range("B5").Value = 69.63786496255987
$ws.Range("B6").Value = 70.53355755666585
$ws.Range("B7").Value = 70.5282024534993
$ws.Range("B8").Value = 64.97578645146216
$ws.Range("B9").Value = 65.18382580505326
$ws.Range("B10").Value = 65.73625264207968
$ws.Range("B11").Value = 66.16906317846573
$ws.Range("B12").Value = 66.23045923029431
$ws.Range("B13").Value = 66.16964904418319
$ws.Range("B14").Value = 67.51521730439242
$ws.Range("B15").Value = 67.7611990824302
$ws.Range("B16").Value = 68.35904731210981
$ws.Range("B17").Value = 68.99338405316506
$ws.Range("B18").Value = 69.37547520870257
$ws.Range("B19").Value = 69.29768626758433
$ws.Range("B20").Value = 69.82827377105221
$ws.Range("B21").Value = 69.96360609809534
$ws.Range("B22").Value = 71.71051551620225
$ws.Range("B23").Value = 70.71721243725725
$ws.Range("B24").Value = 70.21053778326038
$ws.Range("B25").Value = 70.36924478823259
$ws.Range("B26").Value = 70.01125282941254
$ws.Range("B27").Value = 70.23159843081189
$ws.Range("B28").Value = 70.57718893322405

$ws = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws.Range("B2").Value = 71.15485705730194
$ws.Range("C2").Value = 68.69720290577405
$ws.Range("D2").Value = 73.61251120882983
$ws.Range("B3").Value = 71.25759667954611
$ws.Range("C3").Value = 67.78194884683283
$ws.Range("D3").Value = 74.73324451225939
$ws.Range("B4").Value = 71.36033630179028
$ws.Range("C4").Value = 67.10355444391141
$ws.Range("D4").Value = 75.61711815966916
$ws.Range("B5").Value = 71.46307592403446
$ws.Range("C5").Value = 66.54776762097869
$ws.Range("D5").Value = 76.37838422709022
$ws.Range("B6").Value = 71.56581554627863
$ws.Range("C6").Value = 66.07033379827772
$ws.Range("D6").Value = 77.06129729427954

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Range("B2").Value = 0.1048243165158899
$ws.Range("B3").Value = 68.48571157322522
$ws.Range("B4").Value = 68.68777403918652
$ws.Range("B5").Value = 69.63994965683159
$ws.Range("B6").Value = 70.53564225093757
$ws.Range("B7").Value = 70.53028714777102
$ws.Range("B8").Value = 64.97787114573387
$ws.Range("B9").Value = 65.18591049932498
$ws.Range("B10").Value = 65.73833733635139
$ws.Range("B11").Value = 66.17114787273745
$ws.Range("B12").Value = 66.23254392456603
$ws.Range("B13").Value = 66.1717337384549
$ws.Range("B14").Value = 67.51730199866414
$ws.Range("B15").Value = 67.76328377670191
$ws.Range("B16").Value = 68.36113200638152
$ws.Range("B17").Value = 68.99546874743677
$ws.Range("B18").Value = 69.37755990297428
$ws.Range("B19").Value = 69.29977096185604
$ws.Range("B20").Value = 69.83035846532393
$ws.Range("B21").Value = 69.96569079236706
$ws.Range("B22").Value = 71.71260021047397
$ws.Range("B23").Value = 70.71929713152896
$ws.Range("B24").Value = 70.21262247753209
$ws.Range("B25").Value = 70.37132948250431
$ws.Range("B26").Value = 70.01333752368426
$ws.Range("B27").Value = 70.2336831250836
$ws.Range("B28").Value = 70.57927362749577
$ws.Range("B29").Value = 71.15694175157365
$ws.Range("B30").Value = 70.96613420334764
$ws.Range("B31").Value = 70.96345988159651
$ws.Range("B32").Value = 70.96129776131195
$ws.Range("B33").Value = 72.82686752653309

$ws = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws.Range("B2").Value = 71.73526538521779
$ws.Range("C2").Value = 69.36016346818306
$ws.Range("D2").Value = 74.11036730225251
$ws.Range("B3").Value = 71.84008970173367
$ws.Range("C3").Value = 68.48118835864484
$ws.Range("D3").Value = 75.19899104482251
$ws.Range("B4").Value = 71.94491401824956
$ws.Range("C4").Value = 67.8311168247912
$ws.Range("D4").Value = 76.05871121170793
$ws.Range("B5").Value = 72.04973833476545
$ws.Range("C5").Value = 67.29953450069601
$ws.Range("D5").Value = 76.79994216883489
$ws.Range("B6").Value = 72.15456265128134
$ws.Range("C6").Value = 66.84367331130164
$ws.Range("D6").Value = 77.46545199126103
